$d = $word.ActiveDocument

# --- 1. Merge F14 runs into a single run (remove proofErr gramStart/gramEnd split) ---
$f14Xml = @'
<w:p w:rsidR="003474B1" w:rsidRPr="003474B1" w:rsidRDefault="003474B1" w:rsidP="003474B1" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="200" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="003474B1">
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:t xml:space="preserve">F14. The system shall retain information from year to year in order to enable the administrator to contact the previous participants of new call for papers and other events. </w:t>
  </w:r>
</w:p>
'@
$f14Para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "F14.*") { $f14Para = $d.Paragraphs.Item($i); break }
}
$f14Para.Range.InsertXML($f14Xml)

# --- 2. N5 paragraph mark formatting: Times New Roman/sz24 -> Calibri/color000000 ---
$n5Xml = @'
<w:p w:rsidR="003474B1" w:rsidRPr="003474B1" w:rsidRDefault="003474B1" w:rsidP="003474B1" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="200" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="003474B1">
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">N5. The system shall have color codes for the different statuses that an abstract may have. </w:t>
  </w:r>
</w:p>
'@
$n5Para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "N5.*") { $n5Para = $d.Paragraphs.Item($i); break }
}
$n5Para.Range.InsertXML($n5Xml)

# --- 3. Replace trailing br/bookmark paragraph with the new Assumptions section ---
$tailXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve">Assumptions </w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>A1.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> The system will only accept abstract in Word format. This will enable the program to conduct quick and easy word counting and correction mark up.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>A2.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> The system will only accept posters in PowerPoint form. </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>This will enable the program to conduct quick and easy correction mark up.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>A3.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> The system shall save a copy of each entry and correction of the poster’s an abstract so that the administrator may, at her discretion see the progression of a particular submission.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>A4.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> The system shall timestamp each submission and ensure that it is not submitted after the time it is due.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>A5.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> The system shall maintain email addresses as unique identifiers</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> throughout the database. </w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
'@
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertXML($tailXml)

Write-Output "done"
